$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45311

$ws.Range("D20").Value = 286.964
$ws.Range("D21").Value = 318.234
$ws.Range("D22").Value = 410.406
$ws.Range("D23").Value = 615.6079999999999
$ws.Range("D24").Value = 1231.196
$ws.Range("D25").Value = 1478.077
$ws.Range("D26").Value = 246.903
$ws.Range("D27").Value = 298.178
$ws.Range("D28").Value = 277.349
$ws.Range("D29").Value = 328.637
$ws.Range("D30").Value = 277.349
$ws.Range("D31").Value = 328.637
